# feat: add 2022-Q4 data
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# --- 1. Insert a new worksheet for "2022-Q4" right before the existing
#        "2022-Q3" sheet (i.e. right after "总计", ahead of "2022-Q3"),
#        so the tab order becomes 总计, 2022-Q4, 2022-Q3, 2022-Q1. ---
$q3Sheet = $wb.Worksheets.Item(2)
$q4Sheet = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"

# Reuse the header/"A column" cell formatting from the "总计" sheet (bold,
# centered, thin-bordered style) so the new sheet matches the look of the
# other per-quarter sheets.
$summary.Range("B1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q4Sheet.Range("A2:A4").PasteSpecial(-4122)

# Header row
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

# The numeric-looking text columns (fund code / scale / position / ratio /
# market value) must stay TEXT, so force a text number format before
# assigning -- otherwise leading/trailing zeros get silently dropped.
$q4Sheet.Range("B2:B4").NumberFormat = "@"
$q4Sheet.Range("D2:D4").NumberFormat = "@"
$q4Sheet.Range("E2:E4").NumberFormat = "@"
$q4Sheet.Range("F2:F4").NumberFormat = "@"
$q4Sheet.Range("G2:G3").NumberFormat = "@"

# Row 2
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "010714"
$q4Sheet.Range("C2").Value = "东方红远见价值混合A"
$q4Sheet.Range("D2").Value = "16.55"
$q4Sheet.Range("E2").Value = "83.49"
$q4Sheet.Range("F2").Value = "4.41"
$q4Sheet.Range("G2").Value = "0.7299"
$q4Sheet.Range("H2").Value = 6

# Row 3
$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "000928"
$q4Sheet.Range("C3").Value = "中融国企改革灵活配置混合"
$q4Sheet.Range("D3").Value = "0.38"
$q4Sheet.Range("E3").Value = "71.45"
$q4Sheet.Range("F3").Value = "3.69"
$q4Sheet.Range("G3").Value = "0.0140"
$q4Sheet.Range("H3").Value = 9

# Row 4
$q4Sheet.Range("A4").Value = 2
$q4Sheet.Range("B4").Value = "017537"
$q4Sheet.Range("C4").Value = "东方红远见价值混合C"
$q4Sheet.Range("D4").Value = "0.00"
$q4Sheet.Range("E4").Value = "83.49"
$q4Sheet.Range("F4").Value = "4.41"
$q4Sheet.Range("G4").Value = 0
$q4Sheet.Range("H4").Value = 6

# --- 2. Update the "总计" summary sheet: shift the existing two rows down
#        and insert the new 2022-Q4 totals on top (newest quarter first). ---

# Make sure row 4 (new) has the same "A column" look as rows 2/3.
$summary.Range("A3").Copy()
$summary.Range("A4").PasteSpecial(-4122)

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.74

$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.66

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 1.2

# Keep "2022-Q1" as the selected tab, matching the original workbook state.
$wb.Worksheets.Item(4).Select()
